$d = $word.ActiveDocument

# --- Edit 1: text before the "[2]" superscript citation ---
# Before: "The library was modified to incorporate the changes in the decoder architecture"
# After:  "The positional attention layer implementation was missing in the library. We modified
#          the library to add this module in the decoder architecture"
$rng1 = $d.Content
$found1 = $rng1.Find.Execute(
    "The library was modified to incorporate the changes in the decoder architecture",
    $false, $true, $false, $false, $false, $true, 1, $false,
    "The positional attention layer implementation was missing in the library. We modified the library to add this module in the decoder architecture",
    2)
Write-Host "Edit1 found/replaced:" $found1

# --- Edit 2: text right after the "[2]" superscript citation, before the unchanged IWSLT sentence ---
# Before: ". It was however a non-trivial task to make the entire library functional."
# After:  ". To setup and use the library to run the NAR model was a non-trivial task."
$rng2 = $d.Content
$found2 = $rng2.Find.Execute(
    ". It was however a non-trivial task to make the entire library functional.",
    $false, $true, $false, $false, $false, $true, 1, $false,
    ". To setup and use the library to run the NAR model was a non-trivial task.",
    2)
Write-Host "Edit2 found/replaced:" $found2
